$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '[''rTransit'', ''r8'', ''rTransit'', ''r4'', ''r7'', ''r4'', ''rTransit'', ''r8'', ''r9'', ''rTransit'', ''r2'', ''r4'', ''rTransit'', ''rTransit'', ''r11'']'
$ws.Range("B2").Value = '[''rTransit'', ''r8'', ''rTransit'', ''r7'', ''rTransit'', ''r5'', ''r8'', ''r9'', ''r9'', ''rTransit'', ''r3'', ''rTransit'', ''r9'', ''r9'', ''r10'']'
$ws.Range("C2").Value = '[''None'', ''None'', ''None'', ''rTransit'', ''r12'', ''rTransit'', ''r8'', ''r5'', ''rTransit'', ''r9'', ''rTransit'', ''r12'', ''rTransit'', ''r5'', ''r4'']'
$ws.Range("A3").Value = '[''r4'', ''r2'', ''r4'', ''rTransit'', ''r9'', ''rTransit'', ''r7'', ''rTransit'', ''rTransit'', ''r13'', ''rTransit'', ''r10'', ''rTransit'', ''r8'', ''r9'']'
$ws.Range("B3").Value = '[''rTransit'', ''r10'', ''rTransit'', ''r8'', ''r8'', ''rTransit'', ''r2'', ''rTransit'', ''rTransit'', ''r13'', ''rTransit'', ''r9'', ''rTransit'', ''r12'', ''r11'']'
$ws.Range("C3").Value = '[''None'', ''None'', ''None'', ''rTransit'', ''r10'', ''rTransit'', ''r8'', ''rTransit'', ''r11'', ''rTransit'', ''r8'', ''rTransit'', ''r2'', ''rTransit'', ''r9'']'
$ws.Range("A4").Value = '[''rTransit'', ''r10'', ''rTransit'', ''r2'', ''rTransit'', ''r10'', ''r9'', ''rTransit'', ''r3'', ''rTransit'', ''r12'', ''rTransit'', ''r13'', ''rTransit'', ''r8'']'
$ws.Range("B4").Value = '[''rTransit'', ''r8'', ''rTransit'', ''r3'', ''rTransit'', ''r9'', ''r10'', ''r10'', ''rTransit'', ''r13'', ''rTransit'', ''rTransit'', ''r7'', ''r7'', ''r7'']'
$ws.Range("C4").Value = '[''None'', ''None'', ''None'', ''rTransit'', ''rTransit'', ''r11'', ''rTransit'', ''r3'', ''rTransit'', ''r11'', ''rTransit'', ''rTransit'', ''r2'', ''rTransit'', ''r10'']'
$ws.Range("A5").Value = '[''rTransit'', ''r10'', ''rTransit'', ''r11'', ''rTransit'', ''r3'', ''rTransit'', ''r10'', ''r10'', ''rTransit'', ''r4'', ''rTransit'', ''r9'', ''r8'', ''r5'']'
$ws.Range("B5").Value = '[''r11'', ''rTransit'', ''r10'', ''rTransit'', ''r3'', ''rTransit'', ''r12'', ''rTransit'', ''r7'', ''rTransit'', ''r8'', ''r9'', ''rTransit'', ''r11'', ''r12'']'
$ws.Range("C5").Value = '[''None'', ''None'', ''None'', ''r7'', ''r7'', ''rTransit'', ''rTransit'', ''r11'', ''rTransit'', ''rTransit'', ''r7'', ''r2'', ''rTransit'', ''r9'', ''r8'']'
$ws.Range("A6").Value = '[''r1'', ''rTransit'', ''r12'', ''r10'', ''rTransit'', ''r13'', ''r13'', ''rTransit'', ''r9'', ''rTransit'', ''r13'', ''rTransit'', ''r9'', ''rTransit'', ''r13'']'
$ws.Range("B6").Value = '[''r11'', ''rTransit'', ''r10'', ''rTransit'', ''r8'', ''rTransit'', ''r3'', ''rTransit'', ''r12'', ''rTransit'', ''r4'', ''r5'', ''r4'', ''rTransit'', ''r8'']'
$ws.Range("C6").Value = '[''None'', ''None'', ''None'', ''r2'', ''r7'', ''rTransit'', ''r12'', ''rTransit'', ''r4'', ''rTransit'', ''rTransit'', ''r11'', ''r11'', ''r13'', ''r13'']'
$ws.Range("A7").Value = '[''rTransit'', ''r10'', ''rTransit'', ''r1'', ''r2'', ''rTransit'', ''rTransit'', ''r13'', ''rTransit'', ''rTransit'', ''r5'', ''rTransit'', ''r10'', ''rTransit'', ''r11'']'
$ws.Range("B7").Value = '[''rTransit'', ''rTransit'', ''r3'', ''rTransit'', ''r7'', ''r2'', ''rTransit'', ''r12'', ''rTransit'', ''r9'', ''rTransit'', ''r11'', ''rTransit'', ''r3'', ''r3'']'
$ws.Range("C7").Value = '[''None'', ''None'', ''None'', ''rTransit'', ''rTransit'', ''r13'', ''rTransit'', ''r12'', ''rTransit'', ''r2'', ''rTransit'', ''r5'', ''rTransit'', ''r10'', ''r10'']'
$ws.Range("A8").Value = '[''r2'', ''rTransit'', ''r9'', ''rTransit'', ''r3'', ''rTransit'', ''r9'', ''rTransit'', ''r5'', ''rTransit'', ''r7'', ''rTransit'', ''r12'', ''rTransit'', ''r9'']'
$ws.Range("B8").Value = '[''rTransit'', ''rTransit'', ''r5'', ''rTransit'', ''r10'', ''r10'', ''rTransit'', ''r5'', ''r5'', ''r8'', ''rTransit'', ''r3'', ''rTransit'', ''r11'', ''r11'']'
$ws.Range("C8").Value = '[''None'', ''None'', ''None'', ''rTransit'', ''r5'', ''r3'', ''rTransit'', ''r11'', ''r12'', ''r12'', ''rTransit'', ''r8'', ''rTransit'', ''r3'', ''r5'']'
$ws.Range("A9").Value = '[''r7'', ''rTransit'', ''r5'', ''r5'', ''r4'', ''r3'', ''rTransit'', ''r8'', ''r5'', ''rTransit'', ''r12'', ''rTransit'', ''r7'', ''r2'', ''r3'']'
$ws.Range("B9").Value = '[''rTransit'', ''rTransit'', ''r4'', ''r4'', ''rTransit'', ''rTransit'', ''r13'', ''r13'', ''r11'', ''rTransit'', ''r9'', ''r8'', ''rTransit'', ''r13'', ''r11'']'
$ws.Range("C9").Value = '[''None'', ''None'', ''None'', ''rTransit'', ''r9'', ''rTransit'', ''r2'', ''rTransit'', ''rTransit'', ''r11'', ''rTransit'', ''rTransit'', ''r7'', ''r4'', ''r3'']'
$ws.Range("A10").Value = '[''rTransit'', ''r5'', ''r5'', ''rTransit'', ''r7'', ''rTransit'', ''r8'', ''rTransit'', ''r3'', ''r3'', ''rTransit'', ''r12'', ''rTransit'', ''r7'', ''r7'']'
$ws.Range("B10").Value = '[''rTransit'', ''rTransit'', ''r4'', ''r1'', ''r3'', ''r5'', ''rTransit'', ''r9'', ''rTransit'', ''r3'', ''r3'', ''r4'', ''rTransit'', ''rTransit'', ''r13'']'
$ws.Range("C10").Value = '[''None'', ''None'', ''None'', ''rTransit'', ''r12'', ''r12'', ''rTransit'', ''r8'', ''rTransit'', ''r12'', ''r12'', ''rTransit'', ''r9'', ''rTransit'', ''r4'']'
$ws.Range("A11").Value = '[''rTransit'', ''r12'', ''rTransit'', ''r3'', ''rTransit'', ''r11'', ''r12'', ''rTransit'', ''r2'', ''rTransit'', ''r5'', ''r5'', ''rTransit'', ''r9'', ''r10'']'
$ws.Range("B11").Value = '[''rTransit'', ''rTransit'', ''r3'', ''rTransit'', ''r11'', ''rTransit'', ''r10'', ''rTransit'', ''r11'', ''rTransit'', ''r8'', ''rTransit'', ''r12'', ''rTransit'', ''r5'']'
$ws.Range("C11").Value = '[''None'', ''None'', ''None'', ''rTransit'', ''r10'', ''rTransit'', ''r11'', ''r12'', ''rTransit'', ''r4'', ''rTransit'', ''rTransit'', ''r11'', ''rTransit'', ''r5'']'